# All the prod items, monsters and affixes
#
# The stat block (columns C:R) for the affix rows 77-90 on the "Affixes"
# sheet needs to cycle up by two rows (row 79's stats move to row 77,
# row 80's to row 78, ... and the last two rows (77,78) wrap around to
# become the stats for the last two rows (89,90)). Columns A (name) and
# B (description) stay put on their own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 77
$lastRow = 90
$rowCount = $lastRow - $firstRow + 1

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# Snapshot the C:R values for every affected row *before* writing anything,
# since the destination rows overlap the source rows.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write each row's C:R block from the row two below it (wrapping within
# the 77-90 block), leaving columns A/B untouched.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $r + 2
    if ($srcRow -gt $lastRow) {
        $srcRow = $srcRow - $rowCount
    }
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $val = $srcVals[$col]
        if ($val -eq $null) {
            $ws.Range("$col$r").ClearContents()
        } else {
            $ws.Range("$col$r").Value = $val
        }
    }
}
